# service_es.xlsx - Uniandes 2023: pulir redacción CV
#
# Rewrites the two "PCI Registered Reports / Recommender" description cells
# (E2, E3) on the single worksheet ("Hoja1") to the polished Spanish wording,
# and restores the sheet's selection to match the author's saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: "what" text for the Recommender (PCI Registered Reports) entry
$ws.Range("E2").Value = "Asignación de pares, evaluación de propuestas y emisión de recomendaciones de las fases 1 y 2 de \href{https://www.cos.io/initiatives/registered-reports}{reportes registrados}"

# E3: profile note text
$ws.Range("E3").Value = "Perfil como \href{https://rr.peercommunityin.org/public/user_public_page?userId=1996}{recomendador}"

# Restore the window scroll position / active selection seen in the saved file
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 4
$ws.Range("E3").Select()
